$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45183 to 45184 for rows 2 through 34
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
